$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EXPECTED SCORE (D2): 32.38 -> 27.69
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.69"

# SCORE (G2): 32.74 -> 28.00
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "28.00"

# CLUSTER GOOD TO HAVE MATCH column (F2:F10) reshuffled
$ws.Range("F2").Value = "sql : 2"
$ws.Range("F3").Value = "python : 2"
$ws.Range("F4").Value = "hadoop : 6"
$ws.Range("F5").Value = "engineer : 2"
$ws.Range("F6").Value = "data engineer : 2"
$ws.Range("F7").Value = "access : 1"
$ws.Range("F8").Value = "kafka : 3"
$ws.Range("F9").Value = "spark : 3"
$ws.Range("F10").Value = "tools : 4"
